$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '33.709.88'
$ws.Range("E2").Value = '  -0.31%  '

# Row 3
$ws.Range("D3").Value = '1.764.67'
$ws.Range("E3").Value = '  -0.66%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '224.33'
$ws.Range("E5").Value = '  +1.64%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.543'
$ws.Range("E6").Value = '  -1.38%  '

# Row 7
$ws.Range("E7").Value = '  +0.23%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '31.92'
$ws.Range("E8").Value = '  +2.95%  '

# Row 9
$ws.Range("E9").Value = '  +0.95%  '

# Row 10
$ws.Range("E10").Value = '  -2.96%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0937'

# Row 12
$ws.Range("D12").Value = '2.020.00'
$ws.Range("E12").Value = '  -0.58%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '11.24'
$ws.Range("E13").Value = '  +7.45%  '

# Row 14
$ws.Range("D14").Value = '1.774.10'
$ws.Range("E14").Value = '  -0.11%  '

# Row 15
$ws.Range("D15").Value = '33.700.62'
$ws.Range("E15").Value = '  -0.34%  '

# Row 16
$ws.Range("E16").Value = '  -2.34%  '

# Row 17
$ws.Range("E17").Value = '  -1.73%  '

# Row 18
$ws.Range("E18").Value = '  -1.79%  '

# Row 19
$ws.Range("D19").Value = '0.0₃0770'
$ws.Range("E19").Value = '  -0.23%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '237.37'
$ws.Range("E20").Value = '  -2.80%  '

# Row 21
$ws.Range("E21").Value = '  +0.30%  '

# Row 22
$ws.Range("E22").Value = '  -0.09%  '

# Row 23
$ws.Range("E23").Value = '  -0.81%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.05'
$ws.Range("E24").Value = '  -1.48%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '159.23'
$ws.Range("E25").Value = '  +1.21%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '16.10'
$ws.Range("E26").Value = '  -1.44%  '

# Row 27
$ws.Range("E27").Value = '  +0.90%  '

# Row 28
$ws.Range("E28").Value = '  -0.20%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.22'
$ws.Range("E30").Value = '  +2.19%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.0508'
$ws.Range("E31").Value = '  -2.04%  '

# Row 32
$ws.Range("E32").Value = '  -2.97%  '

# Row 33
$ws.Range("E33").Value = '  +0.42%  '

# Row 34
$ws.Range("E34").Value = '  -1.23%  '

# Row 35
$ws.Range("D35").Value = '1.378.35'
$ws.Range("E35").Value = '  -1.00%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.654'
$ws.Range("E36").Value = '  +2.80%  '

# Row 37
$ws.Range("E37").Value = '  -1.31%  '

# Row 38
$ws.Range("E38").Value = '  -0.42%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.22'
$ws.Range("E39").Value = '  +5.68%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.35'
$ws.Range("E40").Value = '  +0.54%  '

# Row 41
$ws.Range("B41").Value = 'InjectiveProtocol'
$ws.Range("C41").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '13.64'
$ws.Range("E41").Value = '  +16.89%  '

# Row 42
$ws.Range("B42").Value = 'ARBITRUM'
$ws.Range("C42").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.905'
$ws.Range("E42").Value = '  -2.46%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '77.54'
$ws.Range("E43").Value = '  -1.56%  '

# Row 44
$ws.Range("E44").Value = '  -1.70%  '

# Row 45
$ws.Range("B45").Value = 'BabyDogeCoin'
$ws.Range("C45").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D45").Value = '0.0₆0139'
$ws.Range("E45").Value = '  +16.78%  '

# Row 46
$ws.Range("B46").Value = 'WEMIXToken'
$ws.Range("C46").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.08'
$ws.Range("E46").Value = '  +4.76%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0499'
$ws.Range("E47").Value = '  +2.36%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '107.61'

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '5.82'
$ws.Range("E49").Value = '  -0.91%  '

# Row 50
$ws.Range("D50").Value = '1.919.56'
$ws.Range("E50").Value = '  +0.08%  '

# Row 51
$ws.Range("E51").Value = '  +0.52%  '

